$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.769.89"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.318.02"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.97"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.79"
$ws.Range("E6").Value = "  -4.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.315.72"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.36"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "661.36"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.859.41"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.37"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.738.40"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.315.34"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.33"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.95"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.886"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.35"
$ws.Range("E23").Value = "  +4.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.81"
$ws.Range("E24").Value = "  -5.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.33"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.83"
$ws.Range("E26").Value = "  -3.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("E27").Value = "  -4.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.21"
$ws.Range("E28").Value = "  -3.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.13"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.23"
$ws.Range("E31").Value = "  +7.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "566.29"
$ws.Range("E32").Value = "  -5.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.89"
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.678.22"
$ws.Range("E36").Value = "  -7.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.49"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("E38").Value = "  -8.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.19"
$ws.Range("E39").Value = "  +5.07%  "
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.61"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.08"
$ws.Range("E42").Value = "  -6.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.32"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0403"
$ws.Range("E46").Value = "  -3.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.64"
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.35"
$ws.Range("E50").Value = "  -4.48%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "127.10"
$ws.Range("E51").Value = "  -3.32%  "
